$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "71.200.75"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +4.92%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.622.17"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +5.47%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.77"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.99"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +4.24%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +2.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.621.00"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +5.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.166"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +16.46%  "
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.346"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +3.96%  "
$ws.Range("E13").Value = "  +1.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.102.14"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +5.48%  "
$ws.Range("E15").Value = "  +5.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000183"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +8.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "70.843.57"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +4.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.613.75"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +5.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "381.02"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +10.12%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.51"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +6.38%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.89"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +6.72%  "
$ws.Range("E22").Value = "  +1.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.11"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.94%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.45"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +7.18%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.84"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +8.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.57"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +8.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.755.87"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +5.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.996"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0955"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +7.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "528.37"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +6.25%  "
$ws.Range("E32").Value = "  +3.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.32"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +6.58%  "
$ws.Range("E34").Value = "  +4.20%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "165.39"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.26%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.13"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +5.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.89"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +9.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.96"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.37"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +5.35%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.03"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +5.75%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.60"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +9.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.335"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +3.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.14"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "153.61"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.66"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +4.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0270"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +7.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.532"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +4.36%  "
$ws.Range("E51").Value = "  +7.63%  "
